# Applies the "minor comment changes & document finished" edit:
#  - reshapes the helper block (H:L) that feeds the bar chart so the raw
#    byte-count row now sits on row 2 (mirroring B:F) and the compression
#    -ratio row moves down to row 3; row 1 (I1:L1) becomes a mirror of the
#    file-name header (B1:E1) instead of the byte counts
#  - applies the "Bytes" number format to the whole raw-data block (B:F
#    rows 2-3) and the percent format to the ratio row (B4:F4 / I3:L3)
#  - repoints the chart series to the new H3 / I2:L2 / I3:L3 cells
#  - drops the stale "_xlchart.v2.*" helper defined names
#  - widens columns B and C:F to fit the now-formatted byte values
#  - moves the active selection to H25

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- sheet data -----------------------------------------------------------

# Row 1 (I1:L1): now mirrors the file-name labels in B1:E1 (was raw bytes)
$ws.Range("I1").Formula = "=B1"
$ws.Range("J1:L1").Formula = "=C1"

# Row 2: H2 mirrors the "raw" row label, I2:L2 now hold the raw byte counts
$ws.Range("H2").Formula = "=A2"
$ws.Range("I2").Formula = "=B2"
$ws.Range("J2:L2").Formula = "=C2"

# Row 3: H3 gets the "compression ratio" label (moved down from H2),
# I3:L3 now hold the compression-ratio formulas (moved down from I2:L2)
$ws.Range("H3").Value = "compression ratio"
$ws.Range("I3").Formula = "=B4"
$ws.Range("J3").Formula = "=C4"
$ws.Range("K3").Formula = "=D4"
$ws.Range("L3").Formula = "=E4"

# --- number formats ---------------------------------------------------

$bytesFmt = 'General\ "Bytes"'
$ws.Range("B2:F3").NumberFormat = $bytesFmt
$ws.Range("H2:L2").NumberFormat = $bytesFmt
$ws.Range("I1:L1").NumberFormat = $bytesFmt

$ws.Range("B4:F4").NumberFormat = "0%"
$ws.Range("I3:L3").NumberFormat = "0%"

# --- column widths (B and C:F now show formatted byte values) -----------

$ws.Columns("B").ColumnWidth = 10.1
$ws.Columns("C:F").ColumnWidth = 10.5

# --- chart: repoint the single bar series to the relocated cells --------

$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart
$ser = $chart.SeriesCollection().Item(1)
$ser.Formula = "=SERIES(Sheet1!`$H`$3,Sheet1!`$I`$2:`$L`$2,Sheet1!`$I`$3:`$L`$3,1)"

# --- drop the stale "_xlchart.v2.*" tracking defined names --------------

for ($i = $wb.Names.Count; $i -ge 1; $i--) {
    $wb.Names.Item($i).Delete()
}

# --- selection moves to H25 ----------------------------------------------

$ws.Range("H25").Select()
